{"js": "// The document contains a single table with one column; each row holds a\n// single-cell statistic. A handful of summary-row values are being updated,\n// and the three \"per-test breakdown\" rows (which currently hold several\n// tab-separated <w:t> runs crammed into one run) are being collapsed down\n// to a single short value each, matching values that used to live in the\n// first three summary rows.\nconst table = context.document.body.tables.getFirst();\n\n// row index (0-based) -> new text for that row's single cell\nconst updates = [\n  [0, \"0M\"],\n  [1, \"0M\"],\n  [2, \"0M\"],\n  [3, \"200\"],\n  [5, \"0.00070\"],\n  [6, \"0.00022\"],\n  [7, \"0.00007\"],\n  [8, \"0.00030\"],\n  [9, \"0.00041\"],\n  [10, \"0.00045\"],\n  [11, \"0.04456\"],\n  [43, \"100\"],\n  [44, \"0.04\"],\n  [45, \"2362\"],\n];\n\nfor (const [rowIndex, newText] of updates) {\n  const cell = table.getCell(rowIndex, 0);\n  const paragraph = cell.body.paragraphs.getFirst();\n  const range = paragraph.getRange();\n  range.insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# The document contains a single table with one column; each row holds a\n# single-cell statistic. A handful of summary-row values are being updated,\n# and the three \"per-test breakdown\" rows (which currently hold several\n# tab-separated runs of text crammed into one paragraph) are being collapsed\n# down to a single short value each, matching values that used to live in\n# the first three summary rows.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Word COM table rows/cells are 1-based.\n$updates = @(\n    @(1, \"0M\"),\n    @(2, \"0M\"),\n    @(3, \"0M\"),\n    @(4, \"200\"),\n    @(6, \"0.00070\"),\n    @(7, \"0.00022\"),\n    @(8, \"0.00007\"),\n    @(9, \"0.00030\"),\n    @(10, \"0.00041\"),\n    @(11, \"0.00045\"),\n    @(12, \"0.04456\"),\n    @(44, \"100\"),\n    @(45, \"0.04\"),\n    @(46, \"2362\")\n)\n\nforeach ($pair in $updates) {\n    $rowIndex = $pair[0]\n    $newText = $pair[1]\n    $t.Cell($rowIndex, 1).Range.Text = $newText\n}\n"}
